# The Excel report should now show the next quotation batch: the date
# advances to 02/11/2021 and the report time becomes 00:03 (ascending
# order fix). Cells C16 (date) and D16 (time) hold these values.
#
# C16 already carries a date-shaped number format (m/d/yyyy), so simply
# assigning a date-looking string to .Value would make Excel parse it
# into a date serial instead of leaving it as literal text - which is
# not what the source file has (it stores the date/time as plain text
# shared strings). To avoid that automatic type coercion (and avoid
# disturbing the existing cell style), we stage the literal text in a
# scratch cell via a formula (forcing a text result), copy it, and use
# PasteSpecial Values so the destination keeps its original style while
# only its text content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")

$scratch.Formula = "=""02/11/2021"""
$scratch.Copy()
$ws.Range("C16").PasteSpecial(-4163)  # xlPasteValues

$scratch.Formula = "=""00:03"""
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)  # xlPasteValues

$scratch.ClearContents()
$excel.CutCopyMode = $false
